$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.997.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +5.54%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.877.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +4.18%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.12%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''281.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +2.99%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.09%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.5265'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +5.51%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.3530'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +1.16%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.07036'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +7.24%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = '''20.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +2.69%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '''0.8143'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -1.96%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '''0.07784'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +0.22%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '''1.876.31'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +4.13%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''5.211'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +3.58%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '''90.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +3.88%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '''1.000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.12%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '''14.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +5.36%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '''0.000008180'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +3.43%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '''1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.06%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").Value = '''27.012.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +5.33%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '''2.109.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +3.68%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''4.760'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +1.35%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +2.27%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''6.226'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +3.33%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''2.383'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +13.98%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''146.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +3.29%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  +4.15%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  +1.39%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''113.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +4.96%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  +1.91%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''4.373'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +4.96%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''0.08882'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.44%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.04897'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +2.73%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.174'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +4.46%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.7426'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +3.81%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''2.877'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.15%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  +9.23%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.404'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +7.13%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.5294'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +3.87%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.01883'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +1.77%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.9831'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +4.37%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''117.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +3.36%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''6.313'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +3.00%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''8.188'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +3.31%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.9998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -0.11%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.4586'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +1.55%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.1366'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.24%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''9.453'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +2.24%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''36.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +2.62%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''1.520'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +2.94%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.05947'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +2.87%  '
$ws.Range("E51").Style = "Normal"
